$d = $word.ActiveDocument

# --------------------------------------------------------------------
# 1. Version bump: "... Typhoon HIL Control Center 2025.1 SP1" -> "...2025.2"
# --------------------------------------------------------------------
$d.Content.Find.Execute("1 SP1", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2", 2) | Out-Null

# --------------------------------------------------------------------
# 2. New bullet point about the VSource/numpy workaround, inserted right
#    after the "Improvements on Examples documentation." bullet and
#    before the "Added Examples:" bullet - same list (numId 45, ilvl 1).
# --------------------------------------------------------------------
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Examples documentation.*") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $insertRange = $target.Range
    $insertRange.Collapse(0)   # wdCollapseEnd
    $insertRange.InsertParagraphAfter()

    $newPara = $target.Next()
    $newPara.Range.Text = "Temporary workaround on internal variable types of the VSource (due to the new numpy version used on THCC 2025.2)"

    # Match the sibling bullets: ListParagraph style, justified, level 2
    # (ilvl 1) of the same num list (numId 45).
    $newPara.Range.ParagraphFormat.Style = "ListParagraph"
    $newPara.Range.ParagraphFormat.Alignment = 3   # wdAlignParagraphJustify
    $newPara.Range.ListFormat.ListLevelNumber = 2
}
